# Rename strain labels RAN12 -> CRD39 and RAN13 -> CRD38 on the two
# per-wavelength data sheets (Sheet2 = 420nm reading, Sheet3 = 600nm
# reading). RAN11 stays untouched because the author is running all
# biological replicates for RAN11 themselves.
#
# Also restore the last-used-cell selection on each sheet, matching
# where the author was working when they saved.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("Sheet2", "Sheet3")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Strain column is L. Rows 18-29 hold the RAN12 replicates, rows
    # 30-41 hold the RAN13 replicates.
    $ws.Range("L18:L29").Value = "CRD39"
    $ws.Range("L30:L41").Value = "CRD38"
}

# Restore the selections that were active on each sheet when saved.
# Sheet3 is updated first, then Sheet2 last so Sheet2 (the tab that was
# active/selected before and after the edit) ends up the active sheet.
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("V11").Select() | Out-Null

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("Q30").Select() | Out-Null
